# Add two new review rows (37 and 38) to Sheet1, mirroring the existing
# "com.hamxa.shaynachim" / "bitcoin" rows, plus their e-mail hyperlinks,
# and move the sheet's selection/scroll position to just past the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of row 33 (same appid/keyword block, and the style
# pattern - s=1 on col A, s=2 on cols C/D, default elsewhere - that the new
# rows are supposed to end up with) into the two new rows.
$ws.Range("A33:G33").Copy()
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A33:G33").Copy()
$ws.Range("A38:G38").PasteSpecial(-4122)

# Row 37
$ws.Cells.Item(37, 1).Value = "com.hamxa.shaynachim"
$ws.Cells.Item(37, 2).Value = "bitcoin"
$ws.Cells.Item(37, 3).Value = "avishaybar12@gmail.com"
$ws.Cells.Item(37, 4).Value = "stefflugar@gmail.com"
$ws.Cells.Item(37, 5).Value = "27/5/2019 15:59"
$ws.Cells.Item(37, 6).Value = "very hard to find good bitcoin app- this is the one!"
$ws.Cells.Item(37, 7).Value = "no"

# Row 38
$ws.Cells.Item(38, 1).Value = "com.hamxa.shaynachim"
$ws.Cells.Item(38, 2).Value = "bitcoin"
$ws.Cells.Item(38, 3).Value = "nitanfriman@gmail.com"
$ws.Cells.Item(38, 4).Value = "ronoren61@gmail.com"
$ws.Cells.Item(38, 5).Value = "27/5/2019 15:59"
$ws.Cells.Item(38, 6).Value = "marvel guide for beginners in this field"
$ws.Cells.Item(38, 7).Value = "no"

# Hyperlinks for the two e-mail columns on the new rows (mirrors rId33-rId36
# in the target workbook).
$ws.Hyperlinks.Add($ws.Cells.Item(37, 3), "mailto:avishaybar12@gmail.com", "", "", "avishaybar12@gmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(37, 4), "mailto:stefflugar@gmail.com", "", "", "stefflugar@gmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(38, 3), "mailto:nitanfriman@gmail.com", "", "", "nitanfriman@gmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(38, 4), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")

# Adding hyperlinks re-styles the touched cells with the default "Hyperlink"
# look; restore the plain centered style the rest of the table uses by
# re-pasting the formats from row 33 over just the e-mail columns.
$ws.Range("C33:D33").Copy()
$ws.Range("C37:D37").PasteSpecial(-4122)
$ws.Range("C33:D33").Copy()
$ws.Range("C38:D38").PasteSpecial(-4122)

# Move the viewport / selection to just below the newly added rows.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F39").Select()
